$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Switch the deck's theme color scheme from the "Integral" (Red Violet)
#    palette to the standard "Office" palette. This is the COM-level effect
#    of picking a different Design/Colors swatch from the Design tab -
#    PowerPoint stores the 12 theme colors (dk1/lt1/dk2/lt2/accent1-6/
#    hlink/folHlink) on the shared theme part, reachable from any slide's
#    ThemeColorScheme.
# ---------------------------------------------------------------------------
$officeColors = @(
    0,        # dk1      000000
    16777215, # lt1      FFFFFF
    6968388,  # dk2      44546A
    15132391, # lt2      E7E6E6
    13998939, # accent1  5B9BD5
    3243501,  # accent2  ED7D31
    10855845, # accent3  A5A5A5
    49407,    # accent4  FFC000
    12874308, # accent5  4472C4
    4697456,  # accent6  70AD47
    12673797, # hlink    0563C1
    7491477   # folHlink 954F72
)

$themeSlide = $p.Slides.Item(2)
$colorScheme = $themeSlide.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Colors($i).RGB = $officeColors[$i - 1]
}

# ---------------------------------------------------------------------------
# 2) Re-style the three summary tables (slides 14-16) that were using the
#    deck's custom table style with the built-in "No Style, No Grid" table
#    style.
# ---------------------------------------------------------------------------
$newTableStyleId = "{2268D67F-80D4-4F96-A426-A92631651E7D}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($shapeIndex = 1; $shapeIndex -le $slide.Shapes.Count; $shapeIndex++) {
        $shape = $slide.Shapes.Item($shapeIndex)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}
